$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 705.7273
$ws.Range("I2").Value = 726.3
$ws.Range("K2").Value = 726.3
$ws.Range("M2").Value = -613.3
# Row 4 (Leve Item ID 5470)
$ws.Range("H4").Value = 433.33334
$ws.Range("I4").Value = 433.33334
$ws.Range("K4").Value = 433.33334
$ws.Range("M4").Value = -319.33334
# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 207
$ws.Range("I9").Value = 208.53847
$ws.Range("J9").Value = 203
$ws.Range("K9").Value = 208.53847
$ws.Range("L9").Value = 203
$ws.Range("M9").Value = -39.53846999999999
$ws.Range("N9").Value = -541
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 3502.1904
$ws.Range("I138").Value = 1023.8461
$ws.Range("J138").Value = 4613.1724
$ws.Range("K138").Value = 3071.5383
$ws.Range("L138").Value = 13839.5172
$ws.Range("M138").Value = 2068.4617
$ws.Range("N138").Value = -24119.5172

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 24850
$ws.Range("I2").Value = 38000
$ws.Range("K2").Value = 38000
$ws.Range("M2").Value = -37887
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 24850
$ws.Range("I116").Value = 38000
$ws.Range("K116").Value = 38000
$ws.Range("M116").Value = -35706
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3687.9773
$ws.Range("I132").Value = 3001.138
$ws.Range("K132").Value = 9003.414000000001
$ws.Range("M132").Value = -6473.414000000001
# Row 138 (Leve Item ID 42350)
$ws.Range("H138").Value = 74755.664
$ws.Range("J138").Value = 74755.664
$ws.Range("L138").Value = 74755.664
$ws.Range("N138").Value = -85035.664

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 24850
$ws.Range("I3").Value = 38000
$ws.Range("K3").Value = 38000
$ws.Range("M3").Value = -37886
# Row 96 (Leve Item ID 19525)
$ws.Range("H96").Value = 29999.5
$ws.Range("I96").Value = 29999.5
$ws.Range("K96").Value = 29999.5
$ws.Range("M96").Value = -27253.5
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 12856.895
$ws.Range("I99").Value = 12856.895
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 12856.895
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -11358.895
$ws.Range("N99").ClearContents()
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1718.3695
$ws.Range("I134").Value = 1197.75
$ws.Range("J134").Value = 5189.1665
$ws.Range("K134").Value = 3593.25
$ws.Range("L134").Value = 15567.4995
$ws.Range("M134").Value = -1058.25
$ws.Range("N134").Value = -20637.4995

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2880.8
$ws.Range("J31").Value = 4979.778
$ws.Range("L31").Value = 4979.778
$ws.Range("N31").Value = -5569.778
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2880.8
$ws.Range("J34").Value = 4979.778
$ws.Range("L34").Value = 4979.778
$ws.Range("N34").Value = -5383.778
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 1913.579
$ws.Range("I58").Value = 690.5333000000001
$ws.Range("J58").Value = 6500
$ws.Range("K58").Value = 690.5333000000001
$ws.Range("L58").Value = 6500
$ws.Range("M58").Value = -487.5333000000001
$ws.Range("N58").Value = -6906
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 11627482
$ws.Range("J99").Value = 6375
$ws.Range("L99").Value = 6375
$ws.Range("N99").Value = -9371
# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 9309.166999999999
$ws.Range("I105").Value = 13123.5
$ws.Range("K105").Value = 13123.5
$ws.Range("M105").Value = -11376.5
# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 250056000
$ws.Range("J107").Value = 4000
$ws.Range("L107").Value = 4000
$ws.Range("N107").Value = -7840
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 11627482
$ws.Range("J126").Value = 6375
$ws.Range("L126").Value = 19125
$ws.Range("N126").Value = -24065
# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 1913.579
$ws.Range("I136").Value = 690.5333000000001
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 2071.5999
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = 478.4000999999998
$ws.Range("N136").Value = -24600

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (Leve Item ID 4847)
$ws.Range("H2").Value = 416937.4
$ws.Range("I2").Value = 1250056.2
$ws.Range("J2").Value = 378
$ws.Range("K2").Value = 7500337.199999999
$ws.Range("L2").Value = 2268
$ws.Range("M2").Value = -7500224.199999999
$ws.Range("N2").Value = -2494
# Row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 37493780
$ws.Range("I4").Value = 29163086
$ws.Range("K4").Value = 87489258
$ws.Range("M4").Value = -87489146
# Row 11 (Leve Item ID 4745)
$ws.Range("H11").Value = 288.7143
$ws.Range("I11").Value = 264.4
$ws.Range("J11").Value = 349.5
$ws.Range("K11").Value = 793.1999999999999
$ws.Range("L11").Value = 1048.5
$ws.Range("M11").Value = -653.1999999999999
$ws.Range("N11").Value = -1328.5
# Row 59 (Leve Item ID 4694)
$ws.Range("H59").Value = 2979.923
$ws.Range("I59").Value = 2875
$ws.Range("J59").Value = 2999
$ws.Range("K59").Value = 8625
$ws.Range("L59").Value = 8997
$ws.Range("M59").Value = -8085
$ws.Range("N59").Value = -10077
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 5962.4707
$ws.Range("I131").Value = 13225.5
$ws.Range("J131").Value = 2000.8182
$ws.Range("K131").Value = 39676.5
$ws.Range("L131").Value = 6002.4546
$ws.Range("M131").Value = -34636.5
$ws.Range("N131").Value = -16082.4546
# Row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 17915.273
$ws.Range("I140").Value = 23158.625
$ws.Range("J140").Value = 3933
$ws.Range("K140").Value = 69475.875
$ws.Range("L140").Value = 11799
$ws.Range("M140").Value = -64295.875
$ws.Range("N140").Value = -22159

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 714.7083
$ws.Range("I2").Value = 743.3182
$ws.Range("K2").Value = 743.3182
$ws.Range("M2").Value = -630.3182
# Row 97 (Leve Item ID 19940)
$ws.Range("H97").Value = 16661.8
$ws.Range("I97").Value = 17637.643
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 17637.643
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -17141.643
$ws.Range("N97").Value = -3992
# Row 123 (Leve Item ID 34150)
$ws.Range("H123").Value = 22436.143
$ws.Range("J123").Value = 22436.143
$ws.Range("L123").Value = 22436.143
$ws.Range("N123").Value = -27336.143
# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 16217.305
$ws.Range("I126").Value = 25556.428
$ws.Range("J126").Value = 12131.4375
$ws.Range("K126").Value = 76669.284
$ws.Range("L126").Value = 36394.3125
$ws.Range("M126").Value = -74199.284
$ws.Range("N126").Value = -41334.3125

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 20126.285
$ws.Range("I40").Value = 28035.688
$ws.Range("K40").Value = 28035.688
$ws.Range("M40").Value = -27899.688
# Row 43 (Leve Item ID 4314)
$ws.Range("H43").Value = 19249.75
$ws.Range("J43").Value = 19249.75
$ws.Range("L43").Value = 19249.75
$ws.Range("N43").Value = -19635.75
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 3424.5
$ws.Range("I46").Value = 849.7143
$ws.Range("J46").Value = 5999.2856
$ws.Range("K46").Value = 849.7143
$ws.Range("L46").Value = 5999.2856
$ws.Range("M46").Value = -661.7143
$ws.Range("N46").Value = -6375.2856
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 825.5
$ws.Range("I55").Value = 825.5
$ws.Range("K55").Value = 825.5
$ws.Range("M55").Value = -652.5
# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 7985.5
$ws.Range("I93").Value = 7985.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 7985.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -6737.5
$ws.Range("N93").ClearContents()
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 417299.34
$ws.Range("I132").Value = 679299.1
$ws.Range("K132").Value = 2037897.3
$ws.Range("M132").Value = -2035367.3

$ws = $wb.Worksheets.Item("WVR")
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 2254.25
$ws.Range("I136").Value = 1708.5
$ws.Range("K136").Value = 5125.5
$ws.Range("M136").Value = -2575.5
# Row 137 (Leve Item ID 42184)
$ws.Range("H137").Value = 49000
$ws.Range("J137").Value = 49000
$ws.Range("L137").Value = 49000
$ws.Range("N137").Value = -59200
